# Apply the table style change recorded in the target commit.
#
# The deck has a single table (slide 6 / shape "Google Shape;127;p18")
# whose a:tblPr/a:tableStyleId was updated from the built-in style
#   {C59FE356-9915-4BC7-B342-69A5EB31371D}
# to
#   {FD98213A-4E4F-429C-855A-D0CFCAEC105C}
#
# PowerPoint does not allow assigning Table.Style as a plain string
# property - it must go through Table.ApplyStyle("{GUID}").
#
# We scan every slide/shape instead of hard-coding slide/shape indices
# so the script keeps working even if shape ordering differs.

$p = $ppt.ActivePresentation

$oldStyleId = "{C59FE356-9915-4BC7-B342-69A5EB31371D}"
$newStyleId = "{FD98213A-4E4F-429C-855A-D0CFCAEC105C}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $t = $sh.Table
            if ($t.Style -eq $oldStyleId) {
                $t.ApplyStyle($newStyleId)
            }
        }
    }
}
